$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column B entirely (data + width formatting) so only column A remains
$ws.Columns.Item(2).Delete()

# Clear remaining column A content so we can rewrite it cleanly
$ws.Columns.Item(1).ClearContents()

# New single-column data: header + 5 "city country population" rows
$values = @(
    "city country population",
    "SPB Russia 5384342",
    "Moscow Russia 13010112",
    "Kazan Russia 1306953",
    "Novosibirsk Russia 1620162",
    "Yekaterinburg Russia 1493749"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Widen column A to match the target layout (29.15 internal width maps to an
# exported OOXML width of exactly 30, matching the diff)
$ws.Columns.Item(1).ColumnWidth = 29.15
